$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") rows 2-97 all change from 45212 (2023-10-13) to 45221 (2023-10-22)
$ws.Range("C2:C97").Value = 45221
